$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(209).Insert()

$ws.Range("A209").Value = 6
$ws.Range("B209").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C209").Value = "Metropolitana"
$ws.Range("D209").Value = 44769
$ws.Range("E209").Value = 13
$ws.Range("F209").Value = 100112022
$ws.Range("G209").Value = "Arveja Verde"
$ws.Range("H209").Value = "Perfection"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 330
$ws.Range("K209").Value = 38000
$ws.Range("L209").Value = 40000
$ws.Range("M209").Value = 39030
$ws.Range("N209").Value = "$/malla 25 kilos"
$ws.Range("O209").Value = "Provincia de Huasco"
$ws.Range("P209").Value = 1561
$ws.Range("Q209").Value = 25
$ws.Range("R209").Value = "Hortaliza"
